$d = $word.ActiveDocument

# 1. First party header block (SVE.ZIA srl ... denominato "):
#    the trailing straight quote is left untouched by excluding it from
#    the search/replace text (Find/Replace auto-corrects straight quotes
#    to smart quotes, which we must avoid here).
$d.Content.Find.Execute(
    "SVE.ZIA srl con sede legale in VI.LE V.VENETO 4 ELLO -  P.I. n.  03329960136 d" + [char]0x2019 + "ora in poi denominato",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ASSOCIAZIONE LA NOSTRA FAMIGLIA con sede legale in VIA DON LUIGI MONZA  1 PONTE LAMBRO -  P.I. n.  00307430132 d" + [char]0x2019 + "ora in poi denominato",
    2)

# 2. Representative name / codice fiscale block. The leading straight
#    quote is likewise excluded from the search/replace text.
$d.Content.Find.Execute(
    ",  rappresentato dal Sig. MARZIA PANZERI nato a LECCO il 01/01/1970, codice fiscale PNZMRZ76P68E507N.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ",  rappresentato dal Sig. CARLA ANDREOTTI nato a LECCO il 01/01/1970, codice fiscale NDRCRL45559E507I.",
    2)

# 3. Article 1 block - company name (with trailing space)
$d.Content.Find.Execute(
    "SVE.ZIA srl ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ASSOCIAZIONE LA NOSTRA FAMIGLIA ",
    2)

# 4. Article 1 block - address line
$d.Content.Find.Execute(
    "con sede legale in VI.LE V.VENETO 4 ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "con sede legale in VIA DON LUIGI MONZA  1 ",
    2)

# 5. Article 1 block - city line
$d.Content.Find.Execute(
    "23848 ELLO ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " PONTE LAMBRO ",
    2)

# 6. Date/time stamp
$d.Content.Find.Execute(
    "Lecco,  2015-05-21 07:00:18",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lecco,  2015-05-25 11:24:21",
    2)

# 7. Final signature line (exact match, no trailing space)
$d.Content.Find.Execute(
    "SVE.ZIA srl",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ASSOCIAZIONE LA NOSTRA FAMIGLIA",
    2)
